$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing client name (row 2)
$ws.Range("A2").Value = "Fahad Ahmed Mohamed"

# Add a new row (row 3) duplicating the client entry with updated status
$ws.Range("A3").Value = "Fahad Ahmed Mohammed"
$ws.Range("C3").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("E3").Value = "Active"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2025-11-16"
$ws.Range("J3").Style = "Normal"
